$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 290.43332
$ws.Range("I33").Value = 301.19232
$ws.Range("K33").Value = 301.19232
$ws.Range("M33").Value = -72.19232
$ws.Range("H100").Value = 2075.5715
$ws.Range("J100").Value = 2772
$ws.Range("L100").Value = 2772
$ws.Range("N100").Value = -3854
$ws.Range("H121").Value = 1391.0834
$ws.Range("I121").Value = 1800
$ws.Range("J121").Value = 1309.3
$ws.Range("K121").Value = 5400
$ws.Range("L121").Value = 3927.9
$ws.Range("M121").Value = -3653
$ws.Range("N121").Value = -7421.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5286.8877
$ws.Range("I32").Value = 4151.384
$ws.Range("K32").Value = 4151.384
$ws.Range("M32").Value = -3864.384
$ws.Range("H39").Value = 4129
$ws.Range("I39").Value = 2172
$ws.Range("J39").Value = 10000
$ws.Range("K39").Value = 2172
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = -1652
$ws.Range("N39").Value = -11040
$ws.Range("H45").Value = 1621.1316
$ws.Range("I45").Value = 1562.5862
$ws.Range("J45").Value = 1809.7778
$ws.Range("K45").Value = 1562.5862
$ws.Range("L45").Value = 1809.7778
$ws.Range("M45").Value = -1185.5862
$ws.Range("N45").Value = -2563.7778
$ws.Range("H61").Value = 6581.5557
$ws.Range("I61").Value = 4660.143
$ws.Range("K61").Value = 4660.143
$ws.Range("M61").Value = -4448.143
$ws.Range("H102").Value = 2316.6667
$ws.Range("I102").Value = 1821.4286
$ws.Range("K102").Value = 1821.4286
$ws.Range("M102").Value = -199.4286
$ws.Range("H110").Value = 1243.1765
$ws.Range("I110").Value = 1217.4286
$ws.Range("K110").Value = 1217.4286
$ws.Range("M110").Value = 827.5714
$ws.Range("H132").Value = 2818.2856
$ws.Range("I132").Value = 1894.6
$ws.Range("J132").Value = 4049.8667
$ws.Range("K132").Value = 5683.799999999999
$ws.Range("L132").Value = 12149.6001
$ws.Range("M132").Value = -3153.799999999999
$ws.Range("N132").Value = -17209.6001
$ws.Range("H136").Value = 6581.5557
$ws.Range("I136").Value = 4660.143
$ws.Range("K136").Value = 13980.429
$ws.Range("M136").Value = -11430.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7044.2256
$ws.Range("I105").Value = 3306.3635
$ws.Range("K105").Value = 3306.3635
$ws.Range("M105").Value = -1559.3635
$ws.Range("H107").Value = 2311.5293
$ws.Range("I107").Value = 2229.4614
$ws.Range("J107").Value = 2578.25
$ws.Range("K107").Value = 2229.4614
$ws.Range("L107").Value = 2578.25
$ws.Range("M107").Value = -309.4614000000001
$ws.Range("N107").Value = -6418.25
$ws.Range("H122").Value = 35000
$ws.Range("J122").Value = 35000
$ws.Range("L122").Value = 35000
$ws.Range("N122").Value = -44800
$ws.Range("H140").Value = 39571.5
$ws.Range("J140").Value = 39571.5
$ws.Range("L140").Value = 39571.5
$ws.Range("N140").Value = -49931.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1477.75
$ws.Range("I16").Value = 1155.5
$ws.Range("J16").Value = 1800
$ws.Range("K16").Value = 1155.5
$ws.Range("L16").Value = 1800
$ws.Range("M16").Value = -868.5
$ws.Range("N16").Value = -2374
$ws.Range("H107").Value = 743.4231
$ws.Range("I107").Value = 872.2
$ws.Range("J107").Value = 314.16666
$ws.Range("K107").Value = 872.2
$ws.Range("L107").Value = 314.16666
$ws.Range("M107").Value = 1047.8
$ws.Range("N107").Value = -4154.16666
$ws.Range("H113").Value = 1477.75
$ws.Range("I113").Value = 1155.5
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1155.5
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 1014.5
$ws.Range("N113").Value = -6140
$ws.Range("H134").Value = 2641.4375
$ws.Range("I134").Value = 1595.3334
$ws.Range("K134").Value = 4786.0002
$ws.Range("M134").Value = -2251.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 328.7143
$ws.Range("I24").Value = 280.2
$ws.Range("J24").Value = 450
$ws.Range("K24").Value = 840.5999999999999
$ws.Range("L24").Value = 1350
$ws.Range("M24").Value = -610.5999999999999
$ws.Range("N24").Value = -1810
$ws.Range("H93").Value = 4392.3076
$ws.Range("J93").Value = 4392.3076
$ws.Range("L93").Value = 13176.9228
$ws.Range("N93").Value = -16920.9228
$ws.Range("H129").Value = 1846.25
$ws.Range("I129").Value = 600
$ws.Range("J129").Value = 2261.6667
$ws.Range("K129").Value = 1800
$ws.Range("L129").Value = 6785.000100000001
$ws.Range("M129").Value = 3200
$ws.Range("N129").Value = -16785.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4270.7144
$ws.Range("I102").Value = 4400.3076
$ws.Range("J102").Value = 2586
$ws.Range("K102").Value = 4400.3076
$ws.Range("L102").Value = 2586
$ws.Range("M102").Value = -2778.3076
$ws.Range("N102").Value = -5830
$ws.Range("H117").Value = 28500
$ws.Range("J117").Value = 28500
$ws.Range("L117").Value = 28500
$ws.Range("N117").Value = -35384
$ws.Range("H122").Value = 6811.4
$ws.Range("I122").Value = 9769
$ws.Range("J122").Value = 2375
$ws.Range("K122").Value = 29307
$ws.Range("L122").Value = 7125
$ws.Range("M122").Value = -26857
$ws.Range("N122").Value = -12025

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1283.1666
$ws.Range("I68").Value = 1279.8
$ws.Range("J68").Value = 1300
$ws.Range("K68").Value = 1279.8
$ws.Range("L68").Value = 1300
$ws.Range("M68").Value = -530.8
$ws.Range("N68").Value = -2798
$ws.Range("H71").Value = 1283.1666
$ws.Range("I71").Value = 1279.8
$ws.Range("J71").Value = 1300
$ws.Range("K71").Value = 6399
$ws.Range("L71").Value = 6500
$ws.Range("M71").Value = -2655
$ws.Range("N71").Value = -13988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 898.19147
$ws.Range("I113").Value = 410.8421
$ws.Range("J113").Value = 1228.8928
$ws.Range("K113").Value = 1232.5263
$ws.Range("L113").Value = 3686.6784
$ws.Range("M113").Value = 937.4737
$ws.Range("N113").Value = -8026.678400000001
$ws.Range("H122").Value = 2619.1282
$ws.Range("I122").Value = 1773.0714
$ws.Range("K122").Value = 5319.2142
$ws.Range("M122").Value = -2869.2142
$ws.Range("H139").Value = 69826.11
$ws.Range("J139").Value = 69826.11
$ws.Range("L139").Value = 69826.11
$ws.Range("N139").Value = -80106.11
